{"js": "// Replace each \"dividend\u00f7divisor=\" expression in the worksheet table with\n// a new expression, one search/replace per pair (all source values are\n// unique in this document, so a plain text search is unambiguous).\nconst replacements = [\n  [\"513\u00f76=\", \"984\u00f74=\"],\n  [\"845\u00f75=\", \"298\u00f74=\"],\n  [\"168\u00f78=\", \"368\u00f77=\"],\n  [\"122\u00f76=\", \"936\u00f74=\"],\n  [\"201\u00f73=\", \"830\u00f77=\"],\n  [\"543\u00f73=\", \"514\u00f75=\"],\n  [\"974\u00f77=\", \"401\u00f79=\"],\n  [\"410\u00f72=\", \"320\u00f77=\"],\n  [\"473\u00f78=\", \"964\u00f73=\"],\n  [\"929\u00f73=\", \"867\u00f78=\"],\n  [\"708\u00f77=\", \"518\u00f72=\"],\n  [\"790\u00f76=\", \"334\u00f74=\"],\n  [\"360\u00f76=\", \"931\u00f75=\"],\n  [\"855\u00f77=\", \"618\u00f72=\"],\n  [\"958\u00f79=\", \"589\u00f76=\"],\n  [\"143\u00f76=\", \"356\u00f76=\"],\n  [\"565\u00f74=\", \"810\u00f78=\"],\n  [\"678\u00f77=\", \"943\u00f72=\"],\n  [\"812\u00f75=\", \"992\u00f74=\"],\n  [\"193\u00f72=\", \"520\u00f79=\"],\n  [\"924\u00f76=\", \"164\u00f73=\"],\n  [\"888\u00f72=\", \"112\u00f79=\"],\n  [\"881\u00f72=\", \"115\u00f78=\"],\n  [\"875\u00f77=\", \"719\u00f72=\"],\n  [\"740\u00f74=\", \"550\u00f77=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each \"dividend\u00f7divisor=\" expression in the worksheet table with a\n# new expression, one Find/Replace per pair (all source values are unique\n# in this document, so a plain text search is unambiguous).\n$replacements = @(\n    @{ Old = \"513\u00f76=\"; New = \"984\u00f74=\" },\n    @{ Old = \"845\u00f75=\"; New = \"298\u00f74=\" },\n    @{ Old = \"168\u00f78=\"; New = \"368\u00f77=\" },\n    @{ Old = \"122\u00f76=\"; New = \"936\u00f74=\" },\n    @{ Old = \"201\u00f73=\"; New = \"830\u00f77=\" },\n    @{ Old = \"543\u00f73=\"; New = \"514\u00f75=\" },\n    @{ Old = \"974\u00f77=\"; New = \"401\u00f79=\" },\n    @{ Old = \"410\u00f72=\"; New = \"320\u00f77=\" },\n    @{ Old = \"473\u00f78=\"; New = \"964\u00f73=\" },\n    @{ Old = \"929\u00f73=\"; New = \"867\u00f78=\" },\n    @{ Old = \"708\u00f77=\"; New = \"518\u00f72=\" },\n    @{ Old = \"790\u00f76=\"; New = \"334\u00f74=\" },\n    @{ Old = \"360\u00f76=\"; New = \"931\u00f75=\" },\n    @{ Old = \"855\u00f77=\"; New = \"618\u00f72=\" },\n    @{ Old = \"958\u00f79=\"; New = \"589\u00f76=\" },\n    @{ Old = \"143\u00f76=\"; New = \"356\u00f76=\" },\n    @{ Old = \"565\u00f74=\"; New = \"810\u00f78=\" },\n    @{ Old = \"678\u00f77=\"; New = \"943\u00f72=\" },\n    @{ Old = \"812\u00f75=\"; New = \"992\u00f74=\" },\n    @{ Old = \"193\u00f72=\"; New = \"520\u00f79=\" },\n    @{ Old = \"924\u00f76=\"; New = \"164\u00f73=\" },\n    @{ Old = \"888\u00f72=\"; New = \"112\u00f79=\" },\n    @{ Old = \"881\u00f72=\"; New = \"115\u00f78=\" },\n    @{ Old = \"875\u00f77=\"; New = \"719\u00f72=\" },\n    @{ Old = \"740\u00f74=\"; New = \"550\u00f77=\" }\n)\n\n$d = $word.ActiveDocument\n\nforeach ($r in $replacements) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $null = $rng.Find.Execute(\n        $r.Old,    # FindText\n        $false,    # MatchCase\n        $false,    # MatchWholeWord\n        $false,    # MatchWildcards\n        $false,    # MatchSoundsLike\n        $false,    # MatchAllWordForms\n        $true,     # Forward\n        1,         # Wrap (wdFindContinue)\n        $false,    # Format\n        $r.New,    # ReplaceWith\n        2          # Replace (wdReplaceAll)\n    )\n}\n"}
